$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.991719663143158
$ws.Range("B1").Value = 3.110146522521973
$ws.Range("C1").Value = 3.958924055099487
$ws.Range("D1").Value = 2.037263631820679
$ws.Range("E1").Value = 1.206307649612427
